$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.064003666814281
$ws.Range("D2").Value = 1.065449079084737
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.074230385536447
$ws.Range("I2").Value = 1.045399440795534
$ws.Range("J2").Value = 1.068966137562483
$ws.Range("K2").Value = 1.068162289781367
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.076920147860174
$ws.Range("N2").Value = 1.070484192089951

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.065766352844172
$ws.Range("D3").Value = 1.066832518187432
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.075750428827662
$ws.Range("I3").Value = 1.045841419034526
$ws.Range("J3").Value = 1.070380011428344
$ws.Range("K3").Value = 1.069359710630731
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.078255551151607
$ws.Range("N3").Value = 1.071900073818876

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.066903966078858
$ws.Range("D4").Value = 1.067724896835661
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.076731410998891
$ws.Range("I4").Value = 1.046124747480966
$ws.Range("J4").Value = 1.071291606865056
$ws.Range("K4").Value = 1.070131194067234
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.079116554021518
$ws.Range("N4").Value = 1.07281296382567

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.06738152430279
$ws.Range("D5").Value = 1.068099393439709
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.07714320796616
$ws.Range("I5").Value = 1.046243225971316
$ws.Range("J5").Value = 1.071674069307135
$ws.Range("K5").Value = 1.07045473903141
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.079477790048546
$ws.Range("N5").Value = 1.073195969408285

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.067461668140125
$ws.Range("D6").Value = 1.068162234709972
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.077212315171677
$ws.Range("I6").Value = 1.046263082047318
$ws.Range("J6").Value = 1.071738241490395
$ws.Range("K6").Value = 1.070509017813594
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.079538400694178
$ws.Range("N6").Value = 1.073260232723406

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.066910349943307
$ws.Range("D7").Value = 1.067729903453785
$ws.Range("E7").Value = 0.9943035907978918
$ws.Range("F7").Value = 1.076736915823047
$ws.Range("I7").Value = 1.046126333076872
$ws.Range("J7").Value = 1.071296720365128
$ws.Range("K7").Value = 1.070135520368143
$ws.Range("L7").Value = 0.9968970624459044
$ws.Range("M7").Value = 1.079121383726954
$ws.Range("N7").Value = 1.072818084587498

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.064599997189373
$ws.Range("D8").Value = 1.065917204817564
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.074744633997703
$ws.Range("I8").Value = 1.045549362975329
$ws.Range("J8").Value = 1.069444648251801
$ws.Range("K8").Value = 1.068567658902249
$ws.Range("L8").Value = 0.9958175282591057
$ws.Range("M8").Value = 1.077372100439338
$ws.Range("N8").Value = 1.070963382319334

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060505502052624
$ws.Range("D9").Value = 1.062701078902071
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.071213642646732
$ws.Range("I9").Value = 1.044512094882534
$ws.Range("J9").Value = 1.066155438276693
$ws.Range("K9").Value = 1.065778945782639
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.07426547163268
$ws.Range("N9").Value = 1.067669501288784

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057759193556481
$ws.Range("D10").Value = 1.060541578458876
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.068845245045325
$ws.Range("I10").Value = 1.04380647888322
$ws.Range("J10").Value = 1.063944651051682
$ws.Range("K10").Value = 1.063901725016347
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.07217745486262
$ws.Range("N10").Value = 1.065455574492332

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056565853204719
$ws.Range("D11").Value = 1.059602681955827
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.067816122642857
$ws.Range("I11").Value = 1.043497532955141
$ws.Range("J11").Value = 1.062982922101724
$ws.Range("K11").Value = 1.063084436638813
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.071269151867582
$ws.Range("N11").Value = 1.064492479776954

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056121948631532
$ws.Range("D12").Value = 1.059253347049175
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.067433307083845
$ws.Range("I12").Value = 1.043382259089059
$ws.Range("J12").Value = 1.062625010894182
$ws.Range("K12").Value = 1.062780179935253
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.070931126711488
$ws.Range("N12").Value = 1.064134060294444

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.056217197144464
$ws.Range("D13").Value = 1.05932830734714
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.067515447606013
$ws.Range("I13").Value = 1.043407009233333
$ws.Range("J13").Value = 1.062701815173957
$ws.Range("K13").Value = 1.062845474977046
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.071003663475741
$ws.Range("N13").Value = 1.064210973645106

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056529173174886
$ws.Range("D14").Value = 1.059573817874073
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.067784490367079
$ws.Range("I14").Value = 1.043488014971309
$ws.Range("J14").Value = 1.062953351052923
$ws.Range("K14").Value = 1.063059300617049
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.07124122375601
$ws.Range("N14").Value = 1.064462866733872

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056721305907839
$ws.Range("D15").Value = 1.059725006838248
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.067950182695599
$ws.Range("I15").Value = 1.043537856532877
$ws.Range("J15").Value = 1.063108239762518
$ws.Range("K15").Value = 1.063190955247796
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.071387507148573
$ws.Range("N15").Value = 1.0646179754032

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057838302202521
$ws.Range("D16").Value = 1.060603808347923
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.068913467649126
$ws.Range("I16").Value = 1.043826910316884
$ws.Range("J16").Value = 1.06400838296717
$ws.Range("K16").Value = 1.063955871162732
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.072237646785381
$ws.Range("N16").Value = 1.065519396914451

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058537835054423
$ws.Range("D17").Value = 1.061154025678338
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.069516739549476
$ws.Range("I17").Value = 1.044007309532587
$ws.Range("J17").Value = 1.064571819065826
$ws.Range("K17").Value = 1.064434485417607
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.072769789027938
$ws.Range("N17").Value = 1.066083633156994

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.058945459567123
$ws.Range("D18").Value = 1.061474590605437
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.069868272370203
$ws.Range("I18").Value = 1.044112204736134
$ws.Range("J18").Value = 1.064900034076742
$ws.Range("K18").Value = 1.064713225525508
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.073079776388807
$ws.Range("N18").Value = 1.066412314270872

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059084381496202
$ws.Range("D19").Value = 1.061583833019761
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.069988077807419
$ws.Range("I19").Value = 1.0441479157131
$ws.Range("J19").Value = 1.065011874803903
$ws.Range("K19").Value = 1.064808196563464
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.073185406246596
$ws.Range("N19").Value = 1.066524313824694

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.058462823425716
$ws.Range("D20").Value = 1.061095030677255
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.069452050020409
$ws.Range("I20").Value = 1.043987988431978
$ws.Range("J20").Value = 1.064511412053053
$ws.Range("K20").Value = 1.064383178944667
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.072712736899468
$ws.Range("N20").Value = 1.066023140359333

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056437321957794
$ws.Range("D21").Value = 1.05950153745268
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.067705279378951
$ws.Range("I21").Value = 1.043464175131845
$ws.Range("J21").Value = 1.062879298958163
$ws.Range("K21").Value = 1.062996353151229
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.071171285976252
$ws.Range("N21").Value = 1.064388709476644

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055160068801041
$ws.Range("D22").Value = 1.058496242481729
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.066603804773692
$ws.Range("I22").Value = 1.043131835592903
$ws.Range("J22").Value = 1.061849169971404
$ws.Range("K22").Value = 1.062120463603186
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.070198399381637
$ws.Range("N22").Value = 1.063357117588566

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.055837525517489
$ws.Range("D23").Value = 1.0590294953503
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.067188026421118
$ws.Range("I23").Value = 1.043308300984521
$ws.Range("J23").Value = 1.06239564064037
$ws.Range("K23").Value = 1.062585166610948
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.070714501610749
$ws.Range("N23").Value = 1.063904364308562

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.058496719191217
$ws.Range("D24").Value = 1.061121689114396
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.069481281501456
$ws.Range("I24").Value = 1.043996719821014
$ws.Range("J24").Value = 1.064538708703441
$ws.Range("K24").Value = 1.064406363452129
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.072738517543224
$ws.Range("N24").Value = 1.066050475774097

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.061566891453701
$ws.Range("D25").Value = 1.063535191541968
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.072128973192924
$ws.Range("I25").Value = 1.044782719339902
$ws.Range("J25").Value = 1.067008891337191
$ws.Range("K25").Value = 1.066503034282016
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.075071543392359
$ws.Range("N25").Value = 1.06852416635052
